# Update capital structure database values for rows 2 and 3 (earnings_debt sheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("earnings_debt")

$ws.Range("D2:D3").Value = 0.371
$ws.Range("F2:F3").Value = 0.113
$ws.Range("G2:G3").Value = 0.2566974088713219
$ws.Range("H2:H3").Value = -0.06971892841458059
$ws.Range("I2:I3").Value = -0.1756001963009653
$ws.Range("J2:J3").Value = -0.1756001963009653
$ws.Range("K2:K3").Value = -175.3
$ws.Range("L2:L3").Value = -0.1924681598594642
$ws.Range("U2:U3").Value = 158.5
$ws.Range("V2:V3").Value = 0.01139648254936079
$ws.Range("W2:W3").Value = -0.8817907444668008
$ws.Range("X2:X3").Value = 0.06216183653903015
$ws.Range("Y2:Y3").Value = -0.943952581005831
$ws.Range("Z2:Z3").Value = 2.537872381753242
$ws.Range("AA2:AA3").Value = -0.4456508884226676
$ws.Range("AB2:AB3").Value = 0.05985155143835302
$ws.Range("AC2:AC3").Value = -0.5055024398610206
$ws.Range("AD2:AD3").Value = 824.2
$ws.Range("AE2:AE3").Value = 88.68329395459614
$ws.Range("AF2:AF3").Value = 912.8832939545962
$ws.Range("AG2:AG3").Value = 754.3832939545962
$ws.Range("AH2:AH3").Value = 0.06159522309790968
$ws.Range("AI2:AI3").Value = 0.7733786970977841
$ws.Range("AJ2:AJ3").Value = 0.05145095234661526
$ws.Range("AK2:AK3").Value = 0.7382284243391444
$ws.Range("AL2:AL3").Value = 29.4
$ws.Range("AM2:AM3").Value = 9.599999999999998
$ws.Range("AN2:AN3").Value = -6.662894098625708
$ws.Range("AO2:AO3").Value = -5.41156462585034
$ws.Range("AP2:AP3").Value = -6.098490654442977
$ws.Range("AQ2:AQ3").Value = -16.57291666666667

$wb.Save()
